# S05/G02: Broker connection flow (Connect Zerodha)
# Update rows 37-39 (S05/G02 tasks) with implementation details now that the
# broker connection flow has been implemented.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37: S05_G02_TB001 - backend OAuth-like flow
$ws.Range("G37").Value2 = "implemented"
$ws.Range("F37").Value2 = "Implemented backend Zerodha OAuth-like flow using JSON config (kite_config.json) and KiteConnect.generate_session."
$ws.Range("H37").Value2 = "Backend exposes login-url, connect, and status endpoints that also call Kite profile to verify the connection."
$ws.Range("I37").Value2 = "Refine error reporting and add retries or token refresh handling as needed."

# Row 38: S05_G02_TB002 - secure token storage
$ws.Range("G38").Value2 = "implemented"
$ws.Range("F38").Value2 = "Access token is encrypted with an env-provided crypto key and stored in broker_connections table."
$ws.Range("H38").Value2 = "BrokerConnection model holds one encrypted access token per broker (currently zerodha)."
$ws.Range("I38").Value2 = "Consider stronger encryption (e.g., cryptography.fernet) for multi-user or hosted deployments."

# Row 39: S05_G02_TF003 - frontend Connect Zerodha UI
$ws.Range("G39").Value2 = "implemented"
$ws.Range("F39").Value2 = "Added Zerodha connection section on Settings page with login button, request_token input, and status chip."
$ws.Range("H39").Value2 = "Settings page shows Zerodha connection status, last-updated time in IST, and user name/id when available."
$ws.Range("I39").Value2 = "Improve UX (e.g., integrating redirect/callback instead of manual token paste) when feasible."
